$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5285
$ws.Range("I3").Value = 5536
$ws.Range("I4").Value = 1267
$ws.Range("I6").Value = 6031
$ws.Range("I7").Value = 18629
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I4").Value = 74
$ws.Range("I7").Value = 596
$ws.Range("I8").Value = 1115
$ws.Range("I9").Value = 87
$ws.Range("I10").Value = 134
$ws.Range("I14").Value = 110
$ws.Range("I18").Value = 133
$ws.Range("I19").Value = 510
$ws.Range("I20").Value = 445
$ws.Range("I21").Value = 86
$ws.Range("I23").Value = 184
$ws.Range("I24").Value = 51
$ws.Range("I25").Value = 91
$ws.Range("I29").Value = 1191
$ws.Range("I31").Value = 179
$ws.Range("I33").Value = 849
$ws.Range("I37").Value = 602
$ws.Range("I40").Value = 30
$ws.Range("I42").Value = 632
$ws.Range("I43").Value = 150
$ws.Range("I44").Value = 138
$ws.Range("I53").Value = 193
$ws.Range("I54").Value = 401
$ws.Range("I57").Value = 73
$ws.Range("I60").Value = 98
$ws.Range("I63").Value = 67
$ws.Range("I65").Value = 425
$ws.Range("I66").Value = 52
$ws.Range("I67").Value = 742
$ws.Range("I73").Value = 166
$ws.Range("I76").Value = 274
$ws.Range("I78").Value = 264
$ws.Range("I83").Value = 392
$ws.Range("I84").Value = 158
$ws.Range("I85").Value = 838
$ws.Range("I86").Value = 112
$ws.Range("I89").Value = 217
$ws.Range("I95").Value = 301
$ws.Range("I96").Value = 196
$ws.Range("I97").Value = 148
$ws.Range("I99").Value = 346
$ws.Range("I101").Value = 18629
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 232
$ws.Range("I7").Value = 838
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 346
$ws.Range("I3").Value = 314
$ws.Range("I6").Value = 359
$ws.Range("I7").Value = 1115
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 193
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 188
$ws.Range("I7").Value = 596
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I4").Value = 26
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 217
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 196
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 110
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 199
$ws.Range("I6").Value = 168
$ws.Range("I7").Value = 602
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 98
$ws.Range("I3").Value = 125
$ws.Range("I7").Value = 346
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 170
$ws.Range("I3").Value = 272
$ws.Range("I4").Value = 48
$ws.Range("I5").Value = 20
$ws.Range("I7").Value = 742
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 52
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 179
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 158
$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 125
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 425
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 147
$ws.Range("I7").Value = 392
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 105
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 301
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 192
$ws.Range("I3").Value = 319
$ws.Range("I7").Value = 849
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 91
$ws.Range("I6").Value = 188
$ws.Range("I7").Value = 401
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 350
$ws.Range("I4").Value = 65
$ws.Range("I6").Value = 321
$ws.Range("I7").Value = 1191
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 182
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 510
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 138
$ws = $wb.Worksheets.Item("River North")
$ws.Range("I3").Value = 63
$ws.Range("I7").Value = 274
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 164
$ws.Range("I3").Value = 217
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 632
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 134
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 65
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 264
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 51
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 65
$ws.Range("I7").Value = 184
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 86
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 128
$ws.Range("I3").Value = 136
$ws.Range("I7").Value = 445
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 133
$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 91
$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 52
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 87
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I2").Value = 54
$ws.Range("I7").Value = 166
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I2").Value = 26
$ws.Range("I7").Value = 148
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 112
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 73
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 98
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 150
$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 30
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 74
